$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 43.21270066666667
$ws.Cells.Item(2, 8).Value2 = 129.638102
$ws.Cells.Item(2, 9).Value2 = 0.1487696778665633
$ws.Cells.Item(2, 10).Value2 = 0.1487696778665633
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 31.94073466666667
$ws.Cells.Item(2, 14).Value2 = 95.822204
$ws.Cells.Item(2, 15).Value2 = 0.3365562672414605
$ws.Cells.Item(2, 16).Value2 = 0.3365562672414606
$ws.Cells.Item(2, 17).Value2 = 1380.24540622409
$ws.Cells.Item(2, 18).Value2 = 12422.20865601681
$ws.Cells.Item(2, 19).Value2 = 0.05006936746148508
$ws.Cells.Item(2, 20).Value2 = 0.05006936746148508

# Row 3
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 43.21270066666667
$ws.Cells.Item(3, 8).Value2 = 129.638102
$ws.Cells.Item(3, 9).Value2 = 0.1487696778665633
$ws.Cells.Item(3, 10).Value2 = 0.1487696778665633
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 38.58528233333333
$ws.Cells.Item(3, 14).Value2 = 115.755847
$ws.Cells.Item(3, 15).Value2 = 0.406569189096231
$ws.Cells.Item(3, 16).Value2 = 0.406569189096231
$ws.Cells.Item(3, 17).Value2 = 1667.374255609155
$ws.Cells.Item(3, 18).Value2 = 15006.36830048239
$ws.Cells.Item(3, 19).Value2 = 0.06048516729231614
$ws.Cells.Item(3, 20).Value2 = 0.06048516729231614

# Row 4
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 43.21270066666667
$ws.Cells.Item(4, 8).Value2 = 129.638102
$ws.Cells.Item(4, 9).Value2 = 0.1487696778665633
$ws.Cells.Item(4, 10).Value2 = 0.1487696778665633
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 24.37857333333333
$ws.Cells.Item(4, 14).Value2 = 73.13571999999999
$ws.Cells.Item(4, 15).Value2 = 0.2568745436623085
$ws.Cells.Item(4, 16).Value2 = 0.2568745436623085
$ws.Cells.Item(4, 17).Value2 = 1053.463992133716
$ws.Cells.Item(4, 18).Value2 = 9481.175929203438
$ws.Cells.Item(4, 19).Value2 = 0.03821514311276208
$ws.Cells.Item(4, 20).Value2 = 0.03821514311276208

# Row 5
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 237.7114053333333
$ws.Cells.Item(5, 8).Value2 = 713.134216
$ws.Cells.Item(5, 9).Value2 = 0.8183762794517323
$ws.Cells.Item(5, 10).Value2 = 0.8183762794517323
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 31.94073466666667
$ws.Cells.Item(5, 14).Value2 = 95.822204
$ws.Cells.Item(5, 15).Value2 = 0.3365562672414605
$ws.Cells.Item(5, 16).Value2 = 0.3365562672414606
$ws.Cells.Item(5, 17).Value2 = 7592.676924992452
$ws.Cells.Item(5, 18).Value2 = 68334.09232493206
$ws.Cells.Item(5, 19).Value2 = 0.2754296658112294
$ws.Cells.Item(5, 20).Value2 = 0.2754296658112295

# Row 6
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 237.7114053333333
$ws.Cells.Item(6, 8).Value2 = 713.134216
$ws.Cells.Item(6, 9).Value2 = 0.8183762794517323
$ws.Cells.Item(6, 10).Value2 = 0.8183762794517323
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 38.58528233333333
$ws.Cells.Item(6, 14).Value2 = 115.755847
$ws.Cells.Item(6, 15).Value2 = 0.406569189096231
$ws.Cells.Item(6, 16).Value2 = 0.406569189096231
$ws.Cells.Item(6, 17).Value2 = 9172.161688640106
$ws.Cells.Item(6, 18).Value2 = 82549.45519776095
$ws.Cells.Item(6, 19).Value2 = 0.3327265803122813
$ws.Cells.Item(6, 20).Value2 = 0.3327265803122814

# Row 7
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 237.7114053333333
$ws.Cells.Item(7, 8).Value2 = 713.134216
$ws.Cells.Item(7, 9).Value2 = 0.8183762794517323
$ws.Cells.Item(7, 10).Value2 = 0.8183762794517323
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 24.37857333333333
$ws.Cells.Item(7, 14).Value2 = 73.13571999999999
$ws.Cells.Item(7, 15).Value2 = 0.2568745436623085
$ws.Cells.Item(7, 16).Value2 = 0.2568745436623085
$ws.Cells.Item(7, 17).Value2 = 5795.064927088391
$ws.Cells.Item(7, 18).Value2 = 52155.58434379552
$ws.Cells.Item(7, 19).Value2 = 0.2102200333282216
$ws.Cells.Item(7, 20).Value2 = 0.2102200333282216

# Row 8
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 9.543019333333334
$ws.Cells.Item(8, 8).Value2 = 28.629058
$ws.Cells.Item(8, 9).Value2 = 0.03285404268170446
$ws.Cells.Item(8, 10).Value2 = 0.03285404268170446
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 31.94073466666667
$ws.Cells.Item(8, 14).Value2 = 95.822204
$ws.Cells.Item(8, 15).Value2 = 0.3365562672414605
$ws.Cells.Item(8, 16).Value2 = 0.3365562672414606
$ws.Cells.Item(8, 17).Value2 = 304.8110484448702
$ws.Cells.Item(8, 18).Value2 = 2743.299436003832
$ws.Cells.Item(8, 19).Value2 = 0.01105723396874608
$ws.Cells.Item(8, 20).Value2 = 0.01105723396874608

# Row 9
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 9.543019333333334
$ws.Cells.Item(9, 8).Value2 = 28.629058
$ws.Cells.Item(9, 9).Value2 = 0.03285404268170446
$ws.Cells.Item(9, 10).Value2 = 0.03285404268170446
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 38.58528233333333
$ws.Cells.Item(9, 14).Value2 = 115.755847
$ws.Cells.Item(9, 15).Value2 = 0.406569189096231
$ws.Cells.Item(9, 16).Value2 = 0.406569189096231
$ws.Cells.Item(9, 17).Value2 = 368.2200952891251
$ws.Cells.Item(9, 18).Value2 = 3313.980857602126
$ws.Cells.Item(9, 19).Value2 = 0.01335744149163355
$ws.Cells.Item(9, 20).Value2 = 0.01335744149163355

# Row 10
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 9.543019333333334
$ws.Cells.Item(10, 8).Value2 = 28.629058
$ws.Cells.Item(10, 9).Value2 = 0.03285404268170446
$ws.Cells.Item(10, 10).Value2 = 0.03285404268170446
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 24.37857333333333
$ws.Cells.Item(10, 14).Value2 = 73.13571999999999
$ws.Cells.Item(10, 15).Value2 = 0.2568745436623085
$ws.Cells.Item(10, 16).Value2 = 0.2568745436623085
$ws.Cells.Item(10, 17).Value2 = 232.6451966390844
$ws.Cells.Item(10, 18).Value2 = 2093.80676975176
$ws.Cells.Item(10, 19).Value2 = 0.00843936722132484
$ws.Cells.Item(10, 20).Value2 = 0.00843936722132484
